# Applies the "rewriting the scan covering all cases" edit:
#  - Employees sheet: column A (ID) becomes plain numbers 1..24 instead of
#    shared-string labels "E1".."E24"; selection/active view moves off this sheet.
#  - Timings sheet: clears the stray C2 value ("Sunday") that had been
#    duplicated into row 2; becomes the active/selected sheet (tab 2, C3).
#  - ProductsList sheet: view no longer scrolled to A6 (topLeftCell cleared).
#  - Workbook: active tab switches from Employees to Timings.

$wb = $excel.ActiveWorkbook

# --- Employees sheet: replace "E1".."E24" text labels in column A with numbers ---
$employees = $wb.Worksheets.Item("Employees")
for ($row = 2; $row -le 25; $row++) {
    $employees.Cells.Item($row, 1).Value = $row - 1
}
$employees.Range("A6").Select()

# --- Timings sheet: clear the extra C2 value duplicated from C3 ("Sunday") ---
$timings = $wb.Worksheets.Item("Timings")
$timings.Range("C2").ClearContents()
$timings.Range("C3").Select()

# --- ProductsList sheet: reset the scrolled view back to the top (keep selection B24) ---
$products = $wb.Worksheets.Item("ProductsList")
$products.Activate()
$products.Range("B24").Select()
$excel.ActiveWindow.ScrollRow = 1

# --- Make Timings the active sheet/tab ---
$timings.Activate()
$timings.Range("C3").Select()
